# "Updated bm & slides for new bug"
# Target slide: #20 ("Bug Metrics") - resize the title textbox, grow the
# bug table, append a new bug row, and add a "Total: 27" textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

# --- 1) Resize/reposition the "Bug Metrics" title textbox (Shape 1) ---
$title = $s.Shapes.Item(1)
$title.Left = 0
$title.Top = 9.134961128234863
$title.Width = 327.6611328125

# --- 2) Append a new row describing the AWS/Bootstrap bug ---
$tableShape = $s.Shapes.Item(2)
$tbl = $tableShape.Table
$null = $tbl.Rows.Add()
$lastRow = $tbl.Rows.Count

$tbl.Cell($lastRow, 1).Shape.TextFrame.TextRange.Text = "6"
$flush = $tbl.Rows.Count
$tbl.Cell($lastRow, 2).Shape.TextFrame.TextRange.Text = "Bootstrap – does not work on AWS"
$flush = $tbl.Rows.Count
$tbl.Cell($lastRow, 3).Shape.TextFrame.TextRange.Text = "10 - Critical"
$flush = $tbl.Rows.Count
$tbl.Cell($lastRow, 4).Shape.TextFrame.TextRange.Text = "Unresolved"
$flush = $tbl.Rows.Count

# --- 3) Set the bug table's graphic-frame height to its final value
#        (must happen after the row is added, otherwise PowerPoint's
#        row-driven autofit recalculates/overrides the height) ---
$tableShape.Height = 315.70220947265625

# --- 4) Add the "Total: 27" textbox next to the title ---
$totalBox = $s.Shapes.AddTextbox(1, 505.20001220703125, 4.006063461303711, 214.8000030517578, 62.47134017944336)
$totalBox.Fill.Visible = 0
$totalBox.Line.Visible = 0

$tf = $totalBox.TextFrame
$tf.WordWrap = 1
$tf.MarginLeft = 0
$tf.MarginRight = 0
$tf.MarginTop = 0
$tf.MarginBottom = 0
$tf.VerticalAnchor = 4
$tf.AutoSize = 0

$tr = $tf.TextRange
$tr.Font.Size = 48
$tr.LanguageID = "en-SG"
$tr.Text = "Total: 27 "

$underlinePart = $tr.Characters(8, 2)
$underlinePart.Font.Underline = 1

$true
